$wb = $excel.ActiveWorkbook

# --- Product Backlog sheet (rId2) ---
$ws2 = $wb.Worksheets.Item("Product Backlog")
$ws2.Range("G2").Value = 13

# --- Sprint Backlog sheet (rId3) ---
$ws3 = $wb.Worksheets.Item("Sprint Backlog")
$ws3.Range("K3").Value = 4
$ws3.Range("L3").Value = "in progress"

$ws3.Range("J5").Value = 3
$ws3.Range("K5").Value = 3
$ws3.Range("L5").Value = "done"

$ws3.Range("J6").Value = 3
$ws3.Range("K6").Value = 3
$ws3.Range("L6").Value = "done"

$ws3.Range("H7").Value = "very high"

$ws3.Range("K8").Value = "?"
$ws3.Range("L8").Value = "done"

# --- column width best-effort for Sprint Backlog col H ---
$ws3.Columns.Item(8).ColumnWidth = 9.140625

# --- View / selection changes ---
$ws2.Range("C18").Select()
$win = $wb.Windows.Item(1)
$win.Zoom = 85

$ws3.Range("K8").Select()

# --- Rename sheet & reorder active tab ---
$ws3.Name = "Sprint 1 Backlog"

# Ensure Sprint 1 Backlog ends up as the active/selected sheet (activeTab=2)
$ws3.Activate()
